$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-23 (Generation 0-21): Fitness -> 7312
$ws.Range("C2:C23").Value = 7312

# Rows 24-66 (Generation 22-64): Fitness -> 7310
$ws.Range("C24:C66").Value = 7310

# Rows 67-252 (Generation 65-250): Fitness -> 7293
$ws.Range("C67:C252").Value = 7293
